$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.1395280839714665
$ws.Range("C2").Value = 14.721198061084568
$ws.Range("D2").Value = 51.186851682670401
$ws.Range("E2").Value = 26.694415955142347
$ws.Range("F2").Value = 33.783052586840576
$ws.Range("G2").Value = 32.77043279487507
$ws.Range("H2").Value = 8.6255161812289778
$ws.Range("I2").Value = 18.46517818359483
$ws.Range("J2").Value = 23.314491922969154
$ws.Range("K2").Value = 20.067866823626705
$ws.Range("L2").Value = 14.862292629074499
$ws.Range("M2").Value = 31.164050925312857
$ws.Range("N2").Value = 16.314961415998098
$ws.Range("O2").Value = 44.300831548197081
$ws.Range("P2").Value = -1.2788377895068237
$ws.Range("Q2").Value = 14.388507418210489
$ws.Range("R2").Value = 37.859936974475218
$ws.Range("T2").Value = 22.971488175051377
$ws.Range("U2").Value = 39.114371317615053
$ws.Range("V2").Value = 32.163210736235108
$ws.Range("W2").Value = 30.323170734411875
$ws.Range("X2").Value = 18.451139645995248
$ws.Range("Y2").Value = 44.322144631527401
$ws.Range("Z2").Value = 68.065145157018407
$ws.Range("AA2").Value = 57.246063182134286
$ws.Range("AB2").Value = 77.218492902431052
$ws.Range("AC2").Value = 47.006348722821805
$ws.Range("AD2").Value = 64.493279474036711
$ws.Range("AE2").Value = 33.414606112900366
$ws.Range("AF2").Value = 43.735452330638239
$ws.Range("AG2").Value = 41.401378455983952
$ws.Range("AH2").Value = 49.473337417291511
$ws.Range("AI2").Value = 69.248251697343903
$ws.Range("AJ2").Value = 41.436509303174475
$ws.Range("AK2").Value = 23.124652411926323
$ws.Range("AL2").Value = 20.225298716298948
$ws.Range("AM2").Value = 47.442058462426679
$ws.Range("AN2").Value = 48.662927939600451
$ws.Range("AO2").Value = 39.255845183885128
$ws.Range("AP2").Value = 42.670108738778261
$ws.Range("AQ2").Value = 17.934304745660988
$ws.Range("AR2").Value = 31.087843417323768
$ws.Range("AS2").Value = 25.448096124069831
$ws.Range("AT2").Value = 66.540441443949177
$ws.Range("AU2").Value = 26.365284863431082
$ws.Range("AV2").Value = 55.897456492912987
$ws.Range("AW2").Value = 58.826276591353597
$ws.Range("AX2").Value = 57.628507147263086
$ws.Range("AY2").Value = 55.386269823811432
$ws.Range("B3").Value = 16.456989414059773
$ws.Range("C3").Value = 29.78277914308093
$ws.Range("D3").Value = 42.642372979346064
$ws.Range("E3").Value = 158.46780849148919
$ws.Range("F3").Value = 24.436285336968204
$ws.Range("G3").Value = 33.077633906215773
$ws.Range("H3").Value = 24.294246671769919
$ws.Range("I3").Value = 11.312065376296873
$ws.Range("J3").Value = 30.319379737464331
$ws.Range("K3").Value = 26.418274292484313
$ws.Range("L3").Value = 25.951867463551022
$ws.Range("M3").Value = 28.925148636198848
$ws.Range("N3").Value = 27.130758080989875
$ws.Range("O3").Value = 68.09162187056971
$ws.Range("P3").Value = 50.562882676008059
$ws.Range("Q3").Value = 12.729995481394626
$ws.Range("R3").Value = 37.575134000392211
$ws.Range("S3").Value = 28.658185871009966
$ws.Range("T3").Value = 33.220217673694599
$ws.Range("U3").Value = 37.303081688208465
$ws.Range("V3").Value = 24.522394808998996
$ws.Range("W3").Value = 22.217184890442606
$ws.Range("X3").Value = 30.078248037511518
$ws.Range("Y3").Value = 69.017484056872078
$ws.Range("Z3").Value = 69.005386068838064
$ws.Range("AA3").Value = 49.550501564511393
$ws.Range("AB3").Value = 53.498335936949267
$ws.Range("AC3").Value = 34.371923798669101
$ws.Range("AD3").Value = 53.118025245424093
$ws.Range("AE3").Value = 46.382649367119633
$ws.Range("AF3").Value = 61.405807917627065
$ws.Range("AG3").Value = 31.720740089520348
$ws.Range("AH3").Value = 26.899032880622215
$ws.Range("AI3").Value = 54.617673762182363
$ws.Range("AJ3").Value = 55.677965070868574
$ws.Range("AK3").Value = 39.767797107896108
$ws.Range("AL3").Value = 45.928798852634586
$ws.Range("AM3").Value = 57.16200398324527
$ws.Range("AN3").Value = 45.27787682358305
$ws.Range("AO3").Value = 85.456723711747046
$ws.Range("AP3").Value = 49.15444766195094
$ws.Range("AQ3").Value = 17.667256728138359
$ws.Range("AR3").Value = 33.871488645165776
$ws.Range("AS3").Value = 42.154340688149254
$ws.Range("AT3").Value = 48.375312778145947
$ws.Range("AU3").Value = 19.88933471683972
$ws.Range("AV3").Value = 44.724897053244561
$ws.Range("AW3").Value = 38.157182029906807
$ws.Range("AX3").Value = 31.06486330231818
$ws.Range("AY3").Value = 66.731041574707803